$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same green fill used for "Done" rows (2-11) to the Duration column (D2:D11),
# which previously had no fill (or, for D6, the white fill).
$ws.Range("D2:D11").Interior.Color = 5296274

# Row 12 ("Add Redux To Cart") moves from "Active" (yellow) to "Done" (green),
# and its status cell value changes accordingly.
$ws.Range("A12:D12").Interior.Color = 5296274
$ws.Range("C12").Value = "Done"

# Row 13 ("Create Order") moves from "Open" (no fill) to "Done" (green),
# and its status cell value changes accordingly.
$ws.Range("A13:D13").Interior.Color = 5296274
$ws.Range("C13").Value = "Done"

# Update the active selection to reflect where the user ended up.
$ws.Range("C13").Select()
